# Translate Sheet1 into a new "Translated_Sheet1" sheet, and update the
# header on Sheet1 from "hun" to "magyar oszlop".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the header cell on the source sheet ---
$ws.Range("A1").Value = "magyar oszlop"

# --- Create the translated sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Translated_Sheet1"

# Copy the Hungarian column across, plus the translated English column.
$ws2.Range("A1").Value = "magyar oszlop"
$ws2.Range("B1").Value = "translated_magyar oszlop"

$ws2.Range("A2").Value = "Isten áldd meg a magyart!"
$ws2.Range("B2").Value = "['God bless the seed.']"

$ws2.Range("A3").Value = "Csak a Fidesz!"
$ws2.Range("B3").Value = "[`"It's only the Fidesz.`"]"

$ws2.Range("A4").Value = "Államháztartási hiány GDP arányosan"
$ws2.Range("B4").Value = "['State deficit in proportion to GDP']"

# --- Copy formatting from the header row: bold, centered, boxed border ---
$headerRange = $ws2.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null
